{"js": "const body = context.document.body;\n\n// Load all paragraphs so we can identify/delete the six \"gibberish\"\n// paragraphs (\"Sqsqs\", \"Dfsf\", \"Fd\", \"Fds\", \"Fd\", \"Sf\") sitting between the\n// first paragraph and the final paragraph (which originally held the \"dsf\"\n// run plus the _GoBack bookmark).\nbody.paragraphs.load(\"items/text\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs.items;\n\n// These are paragraph indices 1..6 (0-based) - delete from the bottom up so\n// earlier indices remain valid while the collection shrinks.\nfor (let i = 6; i >= 1; i--) {\n  paragraphs[i].delete();\n}\nawait context.sync();\n\n// Re-load; the remaining last paragraph still contains the \"dsf\" run\n// followed by the bookmarkStart/bookmarkEnd for \"_GoBack\".\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst remaining = body.paragraphs.items;\nconst lastParagraph = remaining[remaining.length - 1];\n\n// Build a range that spans from the start of that paragraph up to (but not\n// including) the _GoBack bookmark, i.e. exactly the \"dsf\" run, and delete\n// just that - this removes the run while leaving the paragraph mark and the\n// bookmark itself intact.\nconst bookmarkRange = context.document.bookmarks.getByName(\"_GoBack\").getRange();\nconst paragraphStart = lastParagraph.getRange(\"Start\");\nconst textBeforeBookmark = paragraphStart.expandTo(bookmarkRange);\ntextBeforeBookmark.delete();\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the six \"gibberish\" paragraphs (\"Sqsqs\", \"Dfsf\", \"Fd\", \"Fds\", \"Fd\",\n# \"Sf\") that sit between the first paragraph and the final paragraph\n# (whose text was \"dsf\" and which also carries the _GoBack bookmark).\n# These were paragraphs 2 through 7 (1-based COM indexing); delete from the\n# bottom up so earlier indices stay valid as the collection shrinks.\nfor ($i = 7; $i -ge 2; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\n# The remaining last paragraph still holds the \"dsf\" run plus the\n# bookmarkStart/bookmarkEnd for \"_GoBack\". Clear just the run's text while\n# leaving the paragraph mark (and bookmark) in place.\n$last = $d.Paragraphs.Item($d.Paragraphs.Count)\n$r = $last.Range\n$r.MoveEnd(1, -1) | Out-Null\n$r.Text = \"\"\n"}
